$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the date string from 01.16.19 to 01.16.20 across the sheet
# (s2cDNADate column A and libraryDate column D, rows 2-32).
# Force the cells to remain plain text (not get auto-converted to a date
# serial number) while keeping the cells' original (unstyled) appearance.
$rangeA = $ws.Range("A2:A32")
$rangeA.NumberFormat = "@"
$rangeA.Value = "01.16.20"
$rangeA.Style = "Normal"

$rangeD = $ws.Range("D2:D32")
$rangeD.NumberFormat = "@"
$rangeD.Value = "01.16.20"
$rangeD.Style = "Normal"

$ws.Activate()
$ws.Range("A2:A32").Select()
